# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would be misinterpreted as a number by Excel's
# automatic type inference (losing a trailing zero or becoming scientific
# notation). Force these to Text format before assignment so the literal
# string is preserved exactly, matching the source data (inline strings).
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"

$ws.Range("D2").Value = "28.136.96"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "1.802.39"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "315.03"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "0.5265"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "0.07998"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "6.345"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "1.005"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "20.58"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "1.813.16"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "7.337"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "92.88"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "0.00001090"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "0.06615"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D21").Value = "17.34"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "5.974"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").Value = "28.192.51"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "158.17"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "20.55"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.401"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.008.04"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "123.12"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").Value = "0.1095"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "5.550"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "0.07267"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "8.865"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2163"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("D39").Value = "0.02304"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "5.040"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").Value = "0.6186"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "1.160"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "0.6016"
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").Value = "13.13"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "3.767"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "126.36"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "0.06825"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "72.78"
$ws.Range("E51").Value = "  -2.15%  "
